# Insurance sheet ("保險", sheet6): add trailer columns (category, date,
# legislator_name, legislator_id, source_file, index) and relabel/fix the
# existing columns (company, name, owner, property_category) per issue #5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1): format new F:K header cells like the existing
# bold/centered/bordered B1:E1 header cells, then set all header labels. ---
$hdr = $ws.Range("F1:K1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# Dates ("2012-02-13") must stay plain text (matching every other sheet),
# not get auto-parsed into a date serial number - force the G column to
# text format before writing into it.
$ws.Range("G2:G7").NumberFormat = "@"

function Set-Row($r, $idx, $company, $name, $owner) {
    $ws.Range("A$r").Value = $idx
    $ws.Range("B$r").Value = $company
    $ws.Range("C$r").Value = $name
    $ws.Range("D$r").Value = $owner
    $ws.Range("E$r").Value = "insurance"
    $ws.Range("F$r").Value = "normal"
    $ws.Range("G$r").Value = "2012-02-13"
    $ws.Range("H$r").Value = "陳明文"
    $ws.Range("I$r").Value = 828
    $ws.Range("J$r").Value = "tmpf4561"
    $ws.Range("K$r").Value = $idx
}

Set-Row 2 133 "國泰人壽" "得意還本終身" "陳明文"
Set-Row 3 134 "國泰人壽" "新富貴保本投資鏈結型保險第7期" "廖素惠"
Set-Row 4 136 "國泰人壽" "創世變額萬能壽險" "廖素惠"
Set-Row 5 137 "國泰人壽" "創世變額萬能壽險" "廖素惠"
Set-Row 6 138 "富邦人壽" "安泰增額養老壽險" "陳明文"
Set-Row 7 139 "台灣人壽" "台灣人壽富利長紅利率變動型年金保險" "陳〇廷"

Write-Output "sheet6 (保險) updated"
